$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Merge the "DESC" + "A" + "RGAR" hyperlink runs into a single "DESCARGAR"
#    run, keeping the original (orange) run formatting.
# ---------------------------------------------------------------------------
$dup = $d.Content.Duplicate
$dup.Find.Execute("DESCARGAR", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$descStart = $dup.Start
$descEnd = $dup.End

# Re-wrap in a fresh Range - reading .Font straight off a post-Find range is
# unreliable in this host, but a brand-new Range(start,end) reads correctly.
$probe = $d.Range($descStart, $descEnd)
$origColor = $probe.Font.Color
$origSize = $probe.Font.Size

$target = $d.Range($descStart, $descEnd)
$target.Text = "DESCARGAR#TMP#"
$target2 = $d.Range($descStart, $descStart + 14)
$target2.Text = "DESCARGAR"
$target3 = $d.Range($descStart, $descStart + 9)
$target3.Font.Color = $origColor
$target3.Font.Size = $origSize

# ---------------------------------------------------------------------------
# 2) Remove the "_GoBack" bookmark from the image paragraph (it is re-added
#    at the very end of the document, see step 4).
# ---------------------------------------------------------------------------
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

# ---------------------------------------------------------------------------
# 3) The "NOTAS" paragraph's paragraph-mark formatting shrinks from 21 to 20
#    half-points (10.5pt -> 10pt).
# ---------------------------------------------------------------------------
$notasPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "NOTAS:*") {
        $notasPara = $p
        break
    }
}
$markRange = $notasPara.Range
$markRange.Font.Size = 10
$markRange.Font.SizeBi = 10

# ---------------------------------------------------------------------------
# 4) Append three new paragraphs after the "NOTAS" block: a blank spacer
#    paragraph, the new "4- Agregar..." paragraph, and the closing
#    "Al salir..." paragraph (which now carries the relocated _GoBack
#    bookmark). They are inserted before the document's pre-existing empty
#    trailing paragraph, which must stay untouched.
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertPos = $lastPara.Range.Start
$w = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

$pPr20 = @"
<w:pPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:after="240" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="2A2513"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="es-AR"/></w:rPr></w:pPr>
"@

$pPr21 = @"
<w:pPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:after="240" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="2A2513"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:eastAsia="es-AR"/></w:rPr></w:pPr>
"@

$rPr20 = '<w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="2A2513"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="es-AR"/></w:rPr>'

# -- Paragraph A: blank spacer paragraph --------------------------------
$paraA = "<w:p $w>$pPr20</w:p>"

# -- Paragraph B: "4- Agregar un atributo usuario :string, ..." ----------
$paraB = "<w:p $w>" + $pPr20 +
    "<w:r>$rPr20<w:t xml:space=`"preserve`">4- Agregar un atributo </w:t></w:r>" +
    "<w:proofErr w:type=`"gramStart`"/>" +
    "<w:r>$rPr20<w:t>usuario :</w:t></w:r>" +
    "<w:proofErr w:type=`"spellStart`"/>" +
    "<w:r>$rPr20<w:t>string</w:t></w:r>" +
    "<w:proofErr w:type=`"spellEnd`"/>" +
    "<w:proofErr w:type=`"gramEnd`"/>" +
    "<w:r>$rPr20<w:t>, que se inicializará en el constructor y antes de ingresar el numero binario el usuario deberá registrarse como 1, 2 o 3.</w:t></w:r>" +
    "</w:p>"

# -- Paragraph C: "Al salir del programa ..." + relocated _GoBack --------
$paraC = "<w:p $w>" + $pPr21 +
    "<w:r>$rPr20<w:t>Al salir del programa mostrar el acumulado de cada usuario.</w:t></w:r>" +
    "<w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/><w:bookmarkEnd w:id=`"0`"/>" +
    "</w:p>"

# -- trailing dummy paragraph: keeps the real final paragraph untouched --
$paraDummy = "<w:p $w></w:p>"

$xml = $paraA + $paraB + $paraC + $paraDummy

$insertRange = $d.Range($insertPos, $insertPos)
$insertRange.InsertXML($xml) | Out-Null

# Remove the dummy empty paragraph that was needed only to keep the
# original trailing empty paragraph from being merged into.
$dummyIndex = $d.Paragraphs.Count - 1
$dummyPara = $d.Paragraphs.Item($dummyIndex)
$dummyPara.Range.Delete()
